$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: "Combatsystem for Axe" is now fully done (Remain = 0) ---
$ws.Range("D21").Value = 6

# --- Row 22: "Bugfixing" is now fully done (Remain = 0); note text updated ---
$ws.Range("D22").Value = 3
$ws.Range("I22").Value = "Bugfixing siehe Bugs.xlsx Didn't fixed all Bugs"

# --- Row 23 stays "Sprite für Springen" (unchanged content) ---

# --- Row 24: rename task + fill in full data (now finished) ---
$ws.Range("A24").Value = "Axe in the Stone for Cave"
$ws.Range("B24").Value = 9
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("F24").Value = "Sascha"

# --- Row 25: brand-new finished task ---
$ws.Range("A25").Value = "Helth UI for Enemies"
$ws.Range("B25").Value = 6
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("F25").Value = "Cedric"

# --- Update the remembered selection to match the author's last cursor position ---
$ws.Range("E34").Select() | Out-Null
